$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.025763182567854
$ws.Range("D2").Value = 1.036250852028061
$ws.Range("E2").Value = 1.047172102482544
$ws.Range("F2").Value = 1.050967896383689
$ws.Range("I2").Value = 1.033887488249528
$ws.Range("J2").Value = 1.030930388913919
$ws.Range("K2").Value = 1.039045129379645
$ws.Range("L2").Value = 1.049935483952231
$ws.Range("M2").Value = 1.053720702712474
$ws.Range("N2").Value = 1.014342986470656

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.026564937148405
$ws.Range("D3").Value = 1.036853143510253
$ws.Range("E3").Value = 1.047993357873456
$ws.Range("F3").Value = 1.051741229860379
$ws.Range("I3").Value = 1.033994849320018
$ws.Range("J3").Value = 1.031372377749227
$ws.Range("K3").Value = 1.03945753397373
$ws.Range("L3").Value = 1.050568505358021
$ws.Range("M3").Value = 1.054306689496411
$ws.Range("N3").Value = 1.01449054639132

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.027084436247216
$ws.Range("D4").Value = 1.037243353731873
$ws.Range("E4").Value = 1.048526113404673
$ws.Range("F4").Value = 1.052242621523716
$ws.Range("I4").Value = 1.034063194436738
$ws.Range("J4").Value = 1.031658421301696
$ws.Range("K4").Value = 1.039724178176769
$ws.Range("L4").Value = 1.050978808175542
$ws.Range("M4").Value = 1.054686192548976
$ws.Range("N4").Value = 1.01458600739883

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.027303001924941
$ws.Range("D5").Value = 1.03740751286386
$ws.Range("E5").Value = 1.048750404527504
$ws.Range("F5").Value = 1.052453642229037
$ws.Range("I5").Value = 1.034091656767008
$ws.Range("J5").Value = 1.031778684004834
$ws.Range("K5").Value = 1.039836223929975
$ws.Range("L5").Value = 1.051151464269797
$ws.Range("M5").Value = 1.054845812995141
$ws.Range("N5").Value = 1.014626133910708

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.027339709878785
$ws.Range("D6").Value = 1.03743508257785
$ws.Range("E6").Value = 1.04878808273854
$ws.Range("F6").Value = 1.052489087290894
$ws.Range("I6").Value = 1.034096419866556
$ws.Range("J6").Value = 1.031798877189775
$ws.Range("K6").Value = 1.039855033846034
$ws.Range("L6").Value = 1.051180463600271
$ws.Range("M6").Value = 1.054872618456397
$ws.Range("N6").Value = 1.0146328710042

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.027087356073599
$ws.Range("D7").Value = 1.03724554678574
$ws.Range("E7").Value = 1.048529109136161
$ws.Range("F7").Value = 1.052245440269435
$ws.Range("I7").Value = 1.034063575813416
$ws.Range("J7").Value = 1.0316600282202
$ws.Range("K7").Value = 1.039725675541608
$ws.Range("L7").Value = 1.05098111456941
$ws.Range("M7").Value = 1.054688325103396
$ws.Range("N7").Value = 1.014586543592653

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.026033991141871
$ws.Range("D8").Value = 1.036454297190206
$ws.Range("E8").Value = 1.047449369130004
$ws.Range("F8").Value = 1.051229041130159
$ws.Range("I8").Value = 1.033924003871813
$ws.Range("J8").Value = 1.031079750544603
$ws.Range("K8").Value = 1.03918454590685
$ws.Range("L8").Value = 1.050149271412859
$ws.Range("M8").Value = 1.053918670228053
$ws.Range("N8").Value = 1.01439285891779

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.024183346401808
$ws.Range("D9").Value = 1.035063831727843
$ws.Range("E9").Value = 1.045557149465961
$ws.Range("F9").Value = 1.049445716524059
$ws.Range("I9").Value = 1.033669477929485
$ws.Range("J9").Value = 1.030057646779522
$ws.Range("K9").Value = 1.03822946643414
$ws.Range("L9").Value = 1.048688862572526
$ws.Range("M9").Value = 1.05256504759623
$ws.Range("N9").Value = 1.014051428695169

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.022953400398522
$ws.Range("D10").Value = 1.034139536081107
$ws.Range("E10").Value = 1.044302801725321
$ws.Range("F10").Value = 1.048262137802435
$ws.Range("I10").Value = 1.03349406600594
$ws.Range("J10").Value = 1.029376604724999
$ws.Range("K10").Value = 1.037591795069927
$ws.Range("L10").Value = 1.047718995313682
$ws.Range("M10").Value = 1.051664489778365
$ws.Range("N10").Value = 1.013823747542387

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.022421748081789
$ws.Range("D11").Value = 1.033739966148158
$ws.Range("E11").Value = 1.043761373670275
$ws.Range("F11").Value = 1.047750920178171
$ws.Range("I11").Value = 1.033416761052823
$ws.Range("J11").Value = 1.029081808053936
$ws.Range("K11").Value = 1.037315468295719
$ws.Range("L11").Value = 1.047299940999931
$ws.Range("M11").Value = 1.051275000866628
$ws.Range("N11").Value = 1.013725150641671

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.022224409139465
$ws.Range("D12").Value = 1.033591648645987
$ws.Range("E12").Value = 1.043560522694611
$ws.Range("F12").Value = 1.047561225654062
$ws.Range("I12").Value = 1.033387844372168
$ws.Range("J12").Value = 1.028972323661372
$ws.Range("K12").Value = 1.037212797977812
$ws.Range("L12").Value = 1.047144423280155
$ws.Range("M12").Value = 1.051130398026066
$ws.Range("N12").Value = 1.013688526428217

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.022266732629888
$ws.Range("D13").Value = 1.033623458671328
$ws.Range("E13").Value = 1.043603594132234
$ws.Range("F13").Value = 1.047601906938356
$ws.Range("I13").Value = 1.033394056235978
$ws.Range("J13").Value = 1.028995807679902
$ws.Range("K13").Value = 1.037234822455485
$ws.Range("L13").Value = 1.047177776098318
$ws.Range("M13").Value = 1.051161412598957
$ws.Range("N13").Value = 1.013696382480809

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.022405433108625
$ws.Range("D14").Value = 1.033727704114307
$ws.Range("E14").Value = 1.043744765959403
$ws.Range("F14").Value = 1.047735235989966
$ws.Range("I14").Value = 1.033414374911272
$ws.Range("J14").Value = 1.029072757701336
$ws.Range("K14").Value = 1.03730698214405
$ws.Range("L14").Value = 1.047287083037017
$ws.Range("M14").Value = 1.051263046498461
$ws.Range("N14").Value = 1.013722123289013

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.022490909710009
$ws.Range("D15").Value = 1.033791946556766
$ws.Range("E15").Value = 1.043831781034392
$ws.Range("F15").Value = 1.047817410252913
$ws.Range("I15").Value = 1.033426867148232
$ws.Range("J15").Value = 1.029120171338127
$ws.Range("K15").Value = 1.037351438144
$ws.Range("L15").Value = 1.047354448951468
$ws.Range("M15").Value = 1.05132567591799
$ws.Range("N15").Value = 1.013737982941216

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.022988704001172
$ws.Range("D16").Value = 1.034166068215516
$ws.Range("E16").Value = 1.044338770795165
$ws.Range("F16").Value = 1.048296092783004
$ws.Range("I16").Value = 1.033499168089446
$ws.Range("J16").Value = 1.029396171609929
$ws.Range("K16").Value = 1.037610129647239
$ws.Range("L16").Value = 1.047746825781406
$ws.Range("M16").Value = 1.051690348688872
$ws.Range("N16").Value = 1.013830290936226

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.023301205448081
$ws.Range("D17").Value = 1.03440092192376
$ws.Range("E17").Value = 1.04465725218273
$ws.Range("F17").Value = 1.048596701867063
$ws.Range("I17").Value = 1.033544159510167
$ws.Range("J17").Value = 1.029569326747468
$ws.Range("K17").Value = 1.037772344547612
$ws.Range("L17").Value = 1.047993196863384
$ws.Range("M17").Value = 1.051919222317878
$ws.Range("N17").Value = 1.013888191113995

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.023483571133046
$ws.Range("D18").Value = 1.034537971297922
$ws.Range("E18").Value = 1.044843182101264
$ws.Range("F18").Value = 1.048772165392476
$ws.Range("I18").Value = 1.033570271899965
$ws.Range("J18").Value = 1.029670334675147
$ws.Range("K18").Value = 1.037866941313262
$ws.Range("L18").Value = 1.048136988108033
$ws.Range("M18").Value = 1.052052764558516
$ws.Range("N18").Value = 1.013921962362828

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.023545768138895
$ws.Range("D19").Value = 1.034584712248958
$ws.Range("E19").Value = 1.044906607338296
$ws.Range("F19").Value = 1.048832014772887
$ws.Range("I19").Value = 1.033579153419825
$ws.Range("J19").Value = 1.029704777332995
$ws.Range("K19").Value = 1.037899192855043
$ws.Range("L19").Value = 1.048186031935666
$ws.Range("M19").Value = 1.052098306441474
$ws.Range("N19").Value = 1.013933477310806

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.023267667794131
$ws.Range("D20").Value = 1.034375717794274
$ws.Range("E20").Value = 1.044623065049724
$ws.Range("F20").Value = 1.048564436584604
$ws.Range("I20").Value = 1.033539345833579
$ws.Range("J20").Value = 1.029550747851219
$ws.Range("K20").Value = 1.037754942528014
$ws.Range("L20").Value = 1.047966754549466
$ws.Range("M20").Value = 1.051894661766561
$ws.Range("N20").Value = 1.013881979067145

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.022364585384327
$ws.Range("D21").Value = 1.033697003637736
$ws.Range("E21").Value = 1.043703187193616
$ws.Range("F21").Value = 1.047695968521665
$ws.Range("I21").Value = 1.033408397143268
$ws.Range("J21").Value = 1.029050097379542
$ws.Range("K21").Value = 1.037285733737888
$ws.Range("L21").Value = 1.047254891052005
$ws.Range("M21").Value = 1.051233115882651
$ws.Range("N21").Value = 1.013714543284462

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.021797594841898
$ws.Range("D22").Value = 1.033270851905178
$ws.Range("E22").Value = 1.043126326006316
$ws.Range("F22").Value = 1.047151054370684
$ws.Range("I22").Value = 1.033324895109812
$ws.Range("J22").Value = 1.028735413045139
$ws.Range("K22").Value = 1.036990549235298
$ws.Range("L22").Value = 1.04680811186122
$ws.Range("M22").Value = 1.050817585564875
$ws.Range("N22").Value = 1.013609264596078

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.022098089565833
$ws.Range("D23").Value = 1.033496707099901
$ws.Range("E23").Value = 1.043431987885094
$ws.Range("F23").Value = 1.047439816190997
$ws.Range("I23").Value = 1.033369271726041
$ws.Range("J23").Value = 1.028902223761762
$ws.Range("K23").Value = 1.037147048171045
$ws.Range("L23").Value = 1.047044881754259
$ws.Range("M23").Value = 1.051037826603676
$ws.Range("N23").Value = 1.013665075148107

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.023282821741951
$ws.Range("D24").Value = 1.034387106260538
$ws.Range("E24").Value = 1.04463851223368
$ws.Range("F24").Value = 1.048579015496118
$ws.Range("I24").Value = 1.033541521329661
$ws.Range("J24").Value = 1.029559142825641
$ws.Range("K24").Value = 1.037762805814954
$ws.Range("L24").Value = 1.047978702423855
$ws.Range("M24").Value = 1.0519057594864
$ws.Range("N24").Value = 1.013884786026706

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.024661117440128
$ws.Range("D25").Value = 1.035422836162961
$ws.Range("E25").Value = 1.046045085764334
$ws.Range("F25").Value = 1.04990582304159
$ws.Range("I25").Value = 1.03373629168148
$ws.Range("J25").Value = 1.030321827301556
$ws.Range("K25").Value = 1.038476551237248
$ws.Range("L25").Value = 1.049065761956773
$ws.Range("M25").Value = 1.052914671752621
$ws.Range("N25").Value = 1.014139709426674
